$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 28 (2020-05-31): 5 hours, chapter 18 exercises
$ws.Range("B28").Value = 5
$ws.Range("C28").Value = "chapter 18 red, exercises 1 and 2 completed"

# Row 29 (2020-06-01): 0 hours, freeday note with leading space
$ws.Range("B29").Value = 0
$ws.Range("D29").Value = " freeday due to doctors appointment"

# Row 30 (2020-06-02): 0 hours, freeday
$ws.Range("B30").Value = 0
$ws.Range("D30").Value = "freeday"

# Update the view state to match the saved selection/scroll position
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws.Range("D29").Select()
